$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").Value = "Button_Quiz"
$ws.Range("E1").Value = "Button_Worksheet"
$ws.Range("G1").Value = "Button_FC"
$ws.Range("I11").Select() | Out-Null
